$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update closing price for BTC-USD (row 3, column D)
$ws.Range("D3").Value = 92562.07000000001

# Update MACRO_SCORE column (N) for rows 2-6
$ws.Range("N2").Value = 54.84087454262382
$ws.Range("N3").Value = 54.84087454262382
$ws.Range("N4").Value = 54.84087454262382
$ws.Range("N5").Value = 54.84087454262382
$ws.Range("N6").Value = 54.84087454262382
